$wb = $excel.ActiveWorkbook
$wsMeans = $wb.Worksheets.Item("Means")
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# ---------------------------------------------------------------------------
# Sheet "Means": add two new columns (F, G) for the 5-mile and 10-mile radii,
# and update changed values for rows 9 and 10 (Total Cancer Risk / Total
# Respiratory).
# ---------------------------------------------------------------------------

# Header row
$wsMeans.Cells.Item(1, 6).Value2 = "Within 5 miles of HFC production facility"
$wsMeans.Cells.Item(1, 7).Value2 = "Within 10 miles of HFC production facility"

# Row 2 - % White
$wsMeans.Cells.Item(2, 6).Value2 = 80
$wsMeans.Cells.Item(2, 7).Value2 = 80

# Row 3 - % Black or African American
$wsMeans.Cells.Item(3, 6).Value2 = 0.075
$wsMeans.Cells.Item(3, 7).Value2 = 2.6

# Row 4 - % Other
$wsMeans.Cells.Item(4, 6).Value2 = 20
$wsMeans.Cells.Item(4, 7).Value2 = 17

# Row 5 - % Hispanic
$wsMeans.Cells.Item(5, 6).Value2 = 6
$wsMeans.Cells.Item(5, 7).Value2 = 16

# Row 6 - Median Income
$wsMeans.Cells.Item(6, 6).Value2 = 81
$wsMeans.Cells.Item(6, 7).Value2 = 68

# Row 7 - % Below Poverty Line
$wsMeans.Cells.Item(7, 6).Value2 = 6.1
$wsMeans.Cells.Item(7, 7).Value2 = 5.9

# Row 8 - % Below Half the Poverty Line
$wsMeans.Cells.Item(8, 6).Value2 = 4.1
$wsMeans.Cells.Item(8, 7).Value2 = 5.8

# Row 9 - Total Cancer Risk (per million)
$wsMeans.Cells.Item(9, 2).Value2 = 26
$wsMeans.Cells.Item(9, 4).Value2 = 20
$wsMeans.Cells.Item(9, 5).Value2 = 20
$wsMeans.Cells.Item(9, 6).Value2 = 20
$wsMeans.Cells.Item(9, 7).Value2 = 20

# Row 10 - Total Respiratory (hazard quotient)
$wsMeans.Cells.Item(10, 2).Value2 = 0.32
$wsMeans.Cells.Item(10, 3).Value2 = 0.22
$wsMeans.Cells.Item(10, 4).Value2 = 0.2
$wsMeans.Cells.Item(10, 5).Value2 = 0.2
$wsMeans.Cells.Item(10, 6).Value2 = 0.2
$wsMeans.Cells.Item(10, 7).Value2 = 0.21

# ---------------------------------------------------------------------------
# Sheet "Standard Deviations": add two new columns (F, G) for the 5-mile and
# 10-mile radii, and update changed values for rows 9 and 10.
# ---------------------------------------------------------------------------

# Header row
$wsSD.Cells.Item(1, 6).Value2 = "Within 5 mile of HFC production facility SD"
$wsSD.Cells.Item(1, 7).Value2 = "Within 10 mile of HFC production facility SD"

# Row 2 - % White
$wsSD.Cells.Item(2, 6).Value2 = 28
$wsSD.Cells.Item(2, 7).Value2 = 21

# Row 3 - % Black or African American
$wsSD.Cells.Item(3, 6).Value2 = 0.21
$wsSD.Cells.Item(3, 7).Value2 = 4.7

# Row 4 - % Other
$wsSD.Cells.Item(4, 6).Value2 = 28
$wsSD.Cells.Item(4, 7).Value2 = 20

# Row 5 - % Hispanic
$wsSD.Cells.Item(5, 6).Value2 = 13
$wsSD.Cells.Item(5, 7).Value2 = 16

# Row 6 - Median Income
$wsSD.Cells.Item(6, 6).Value2 = 23
$wsSD.Cells.Item(6, 7).Value2 = 20

# Row 7 - % Below Poverty Line
$wsSD.Cells.Item(7, 6).Value2 = 5.6
$wsSD.Cells.Item(7, 7).Value2 = 7.6

# Row 8 - % Below Half the Poverty Line
$wsSD.Cells.Item(8, 6).Value2 = 4.9
$wsSD.Cells.Item(8, 7).Value2 = 6.6

# Row 9 - Total Cancer Risk (per million)
$wsSD.Cells.Item(9, 2).Value2 = 8.6
$wsSD.Cells.Item(9, 3).Value2 = 0.37
$wsSD.Cells.Item(9, 4).Value2 = 0
$wsSD.Cells.Item(9, 5).Value2 = 0
$wsSD.Cells.Item(9, 6).Value2 = 0
$wsSD.Cells.Item(9, 7).Value2 = 0

# Row 10 - Total Respiratory (hazard quotient)
$wsSD.Cells.Item(10, 3).Value2 = 0.034
$wsSD.Cells.Item(10, 4).Value2 = 0.000000000000000028
$wsSD.Cells.Item(10, 5).Value2 = 0.000000000000000016
$wsSD.Cells.Item(10, 6).Value2 = 0.000000000000000013
$wsSD.Cells.Item(10, 7).Value2 = 0.036
